# Auto-generated edit script applying the Lamia_Profits.xlsx diff
# to the corresponding sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 7311.6665
$ws.Range("I18").Value = 985
$ws.Range("J18").Value = 10475
$ws.Range("K18").Value = 985
$ws.Range("L18").Value = 10475
$ws.Range("M18").Value = -701
$ws.Range("N18").Value = -11043
$ws.Range("H92").Value = 1938.4286
$ws.Range("I92").Value = 678.1667
$ws.Range("K92").Value = 678.1667
$ws.Range("M92").Value = 569.8333
$ws.Range("H98").Value = 3264.9744
$ws.Range("I98").Value = 724.9355
$ws.Range("J98").Value = 13107.625
$ws.Range("K98").Value = 724.9355
$ws.Range("L98").Value = 13107.625
$ws.Range("M98").Value = 773.0645
$ws.Range("N98").Value = -16103.625
$ws.Range("H122").Value = 3264.9744
$ws.Range("I122").Value = 724.9355
$ws.Range("J122").Value = 13107.625
$ws.Range("K122").Value = 2174.8065
$ws.Range("L122").Value = 39322.875
$ws.Range("M122").Value = 275.1934999999999
$ws.Range("N122").Value = -44222.875
$ws.Range("H132").Value = 1500.5385
$ws.Range("I132").Value = 1507.5745
$ws.Range("K132").Value = 4522.7235
$ws.Range("M132").Value = -1992.7235

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18184520
$ws.Range("I32").Value = 18520734
$ws.Range("K32").Value = 18520734
$ws.Range("M32").Value = -18520447
$ws.Range("H74").Value = 16669073
$ws.Range("I74").Value = 19609720
$ws.Range("J74").Value = 5411
$ws.Range("K74").Value = 19609720
$ws.Range("L74").Value = 5411
$ws.Range("M74").Value = -19608846
$ws.Range("N74").Value = -7159
$ws.Range("H77").Value = 16669073
$ws.Range("I77").Value = 19609720
$ws.Range("J77").Value = 5411
$ws.Range("K77").Value = 98048600
$ws.Range("L77").Value = 27055
$ws.Range("M77").Value = -98044232
$ws.Range("N77").Value = -35791
$ws.Range("H88").Value = 3129.0527
$ws.Range("I88").Value = 3087.2856
$ws.Range("K88").Value = 3087.2856
$ws.Range("M88").Value = -2681.2856
$ws.Range("H91").Value = 3129.0527
$ws.Range("I91").Value = 3087.2856
$ws.Range("K91").Value = 3087.2856
$ws.Range("M91").Value = -1683.2856
$ws.Range("H110").Value = 9949
$ws.Range("J110").Value = 10827.571
$ws.Range("L110").Value = 10827.571
$ws.Range("N110").Value = -14917.571
$ws.Range("H122").Value = 38462656
$ws.Range("I122").Value = 1108.3334
$ws.Range("J122").Value = 500001250
$ws.Range("K122").Value = 3325.0002
$ws.Range("L122").Value = 1500003750
$ws.Range("M122").Value = -875.0001999999999
$ws.Range("N122").Value = -1500008650

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 5680.778
$ws.Range("I107").Value = 4883
$ws.Range("K107").Value = 4883
$ws.Range("M107").Value = -2963
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = ""
$ws.Range("H109").Value = 81250
$ws.Range("J109").Value = 81250
$ws.Range("L109").Value = 81250
$ws.Range("N109").Value = -84024

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 656.2222
$ws.Range("J7").Value = 224.57143
$ws.Range("L7").Value = 224.57143
$ws.Range("N7").Value = -450.57143
$ws.Range("H31").Value = 119004.555
$ws.Range("I31").Value = 4753.25
$ws.Range("J31").Value = 210405.6
$ws.Range("K31").Value = 4753.25
$ws.Range("L31").Value = 210405.6
$ws.Range("M31").Value = -4458.25
$ws.Range("N31").Value = -210995.6
$ws.Range("H34").Value = 119004.555
$ws.Range("I34").Value = 4753.25
$ws.Range("J34").Value = 210405.6
$ws.Range("K34").Value = 4753.25
$ws.Range("L34").Value = 210405.6
$ws.Range("M34").Value = -4551.25
$ws.Range("N34").Value = -210809.6
$ws.Range("H58").Value = 5380.353
$ws.Range("I58").Value = 2196.6667
$ws.Range("J58").Value = 8962
$ws.Range("K58").Value = 2196.6667
$ws.Range("L58").Value = 8962
$ws.Range("M58").Value = -1993.6667
$ws.Range("N58").Value = -9368
$ws.Range("H99").Value = 3795
$ws.Range("J99").Value = 3795
$ws.Range("L99").Value = 3795
$ws.Range("N99").Value = -6791
$ws.Range("H104").Value = 45000
$ws.Range("J104").Value = 45000
$ws.Range("L104").Value = 45000
$ws.Range("N104").Value = -50242
$ws.Range("H126").Value = 3795
$ws.Range("J126").Value = 3795
$ws.Range("L126").Value = 11385
$ws.Range("N126").Value = -16325
$ws.Range("H134").Value = 6612.375
$ws.Range("I134").Value = 5000
$ws.Range("J134").Value = 6842.7144
$ws.Range("K134").Value = 15000
$ws.Range("L134").Value = 20528.1432
$ws.Range("M134").Value = -12465
$ws.Range("N134").Value = -25598.1432
$ws.Range("H136").Value = 5380.353
$ws.Range("I136").Value = 2196.6667
$ws.Range("J136").Value = 8962
$ws.Range("K136").Value = 6590.000100000001
$ws.Range("L136").Value = 26886
$ws.Range("M136").Value = -4040.000100000001
$ws.Range("N136").Value = -31986

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6957504
$ws.Range("J4").Value = 10785107
$ws.Range("L4").Value = 32355321
$ws.Range("N4").Value = -32355545
$ws.Range("H121").Value = 1890.5
$ws.Range("I121").Value = 1746.5
$ws.Range("K121").Value = 5239.5
$ws.Range("M121").Value = -3929.5
$ws.Range("H122").Value = 6019443.5
$ws.Range("J122").Value = 6213081
$ws.Range("L122").Value = 55917729
$ws.Range("N122").Value = -55922629

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9035.645500000001
$ws.Range("I70").Value = 6685.1113
$ws.Range("J70").Value = 12290.23
$ws.Range("K70").Value = 6685.1113
$ws.Range("L70").Value = 12290.23
$ws.Range("M70").Value = -6415.1113
$ws.Range("N70").Value = -12830.23
$ws.Range("H73").Value = 9035.645500000001
$ws.Range("I73").Value = 6685.1113
$ws.Range("J73").Value = 12290.23
$ws.Range("K73").Value = 6685.1113
$ws.Range("L73").Value = 12290.23
$ws.Range("M73").Value = -5749.1113
$ws.Range("N73").Value = -14162.23
$ws.Range("H136").Value = 48959.4
$ws.Range("J136").Value = 48959.4
$ws.Range("L136").Value = 146878.2
$ws.Range("N136").Value = -151978.2

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 608.5454999999999
$ws.Range("J46").Value = 618
$ws.Range("L46").Value = 618
$ws.Range("N46").Value = -994
$ws.Range("H122").Value = 8786.25
$ws.Range("J122").Value = 13001
$ws.Range("L122").Value = 39003
$ws.Range("N122").Value = -43903
$ws.Range("H132").Value = 13785.714
$ws.Range("I132").Value = 15000
$ws.Range("J132").Value = 13583.333
$ws.Range("K132").Value = 45000
$ws.Range("L132").Value = 40749.999
$ws.Range("M132").Value = -42470
$ws.Range("N132").Value = -45809.999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 20179.055
$ws.Range("J48").Value = 21001
$ws.Range("L48").Value = 21001
$ws.Range("N48").Value = -22139
$ws.Range("H49").Value = 29000
$ws.Range("J49").Value = 29000
$ws.Range("L49").Value = 29000
$ws.Range("N49").Value = -29460
$ws.Range("H136").Value = 1943.2368
$ws.Range("I136").Value = 1509.9706
$ws.Range("K136").Value = 4529.9118
$ws.Range("M136").Value = -1979.9118
